$d = $word.ActiveDocument

# Locate the anchor paragraph: the dashed-line separator paragraph that
# precedes the insertion point (end of the "Observation : same mistake..." block).
$anchor = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -eq "----------------------------------------------------------------------------------------`r") {
        $anchor = $p
    }
}
if ($anchor -eq $null) { throw "Anchor paragraph (dashed separator) not found" }

# Grab a formatted-text template from the anchor run: it carries the exact
# boilerplate rPr/pPr (rFonts cs=Arial, color, szCs, lang bidi, ligatures) used
# throughout this document, so pasting it seeds new paragraphs/runs with the
# right formatting before we overwrite the text.
$srcRange = $d.Range($anchor.Range.Start, $anchor.Range.End - 1)
$template = $srcRange.FormattedText
$templateLen = $srcRange.End - $srcRange.Start

$pos = $anchor.Range.End

# ---- New paragraph 1 of 7 ----
$r = $d.Range($pos, $pos)
$r.InsertParagraphAfter()
$paraStart = $pos
$pos = $paraStart + 1

# ---- New paragraph 2 of 7 ----
$r = $d.Range($pos, $pos)
$r.InsertParagraphAfter()
$paraStart = $pos
$dest = $d.Range($paraStart, $paraStart)
$dest.FormattedText = $template
$pasteEnd = $paraStart + $templateLen
$full = $d.Range($paraStart, $pasteEnd)
$full.Text = "(1)[P35] 5.3.8.1(1)-  CandA(gm)#si | upa# | (GS-5.3-20)"
$paraTextEnd = $full.End
$p = $paraStart
$p = $p + 55  # run 1: "(1)[P35] 5.3.8.1(1)-  CandA(gm)#si | upa# | (GS-5.3-20)" (default formatting)
$pos = $paraTextEnd + 1

# ---- New paragraph 3 of 7 ----
$r = $d.Range($pos, $pos)
$r.InsertParagraphAfter()
$paraStart = $pos
$dest = $d.Range($paraStart, $paraStart)
$dest.FormattedText = $template
$pasteEnd = $paraStart + $templateLen
$full = $d.Range($paraStart, $pasteEnd)
$full.Text = "CandA(gm)#siqdupopaq CandA(gm)#siq CandA(gm)#siqdupa# | "
$paraTextEnd = $full.End
$p = $paraStart
$p = $p + 10  # run 1: "CandA(gm)#" (default formatting)
$sub = $d.Range($p, $p + 6)
$sub.Font.Color = 255  # w:color FF0000
$p = $p + 6
$p = $p + 28  # run 3: "opaq CandA(gm)#siq CandA(gm)" (default formatting)
$sub = $d.Range($p, $p + 6)
$sub.Font.Color = 255  # w:color FF0000
$sub.Select()
$sel = $word.Selection
$find = $sel.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Replacement.Highlight = $true
$find.Execute("#siqdu", $false, $false, $false, $false, $false, $true, 1, $false, "#siqdu", 2) | Out-Null
$p = $p + 6
$p = $p + 6  # run 5: "pa# | " (default formatting)
$pos = $paraTextEnd + 1

# ---- New paragraph 4 of 7 ----
$r = $d.Range($pos, $pos)
$r.InsertParagraphAfter()
$paraStart = $pos
$dest = $d.Range($paraStart, $paraStart)
$dest.FormattedText = $template
$pasteEnd = $paraStart + $templateLen
$full = $d.Range($paraStart, $pasteEnd)
$full.Text = "Observation : same system error as above in Sandhi"
$paraTextEnd = $full.End
$p = $paraStart
$p = $p + 50  # run 1: "Observation : same system error as above in Sandhi" (default formatting)
$pos = $paraTextEnd + 1

# ---- New paragraph 5 of 7 ----
$r = $d.Range($pos, $pos)
$r.InsertParagraphAfter()
$paraStart = $pos
$pos = $paraStart + 1

# ---- New paragraph 6 of 7 ----
$r = $d.Range($pos, $pos)
$r.InsertParagraphAfter()
$paraStart = $pos
$dest = $d.Range($paraStart, $paraStart)
$dest.FormattedText = $template
$pasteEnd = $paraStart + $templateLen
$full = $d.Range($paraStart, $pasteEnd)
$full.Text = "(52)[P36] 5.3.8.2(44)-  CandA(gm)#si | upa# |"
$paraTextEnd = $full.End
$p = $paraStart
$p = $p + 45  # run 1: "(52)[P36] 5.3.8.2(44)-  CandA(gm)#si | upa# |" (default formatting)
$pos = $paraTextEnd + 1

# ---- New paragraph 7 of 7 ----
$r = $d.Range($pos, $pos)
$r.InsertParagraphAfter()
$paraStart = $pos
$dest = $d.Range($paraStart, $paraStart)
$dest.FormattedText = $template
$pasteEnd = $paraStart + $templateLen
$full = $d.Range($paraStart, $pasteEnd)
$full.Text = "CandA(gm)#siqdupopaq CandA(gm)#siq CandA(gm)#siqdupa# |"
$paraTextEnd = $full.End
$p = $paraStart
$p = $p + 10  # run 1: "CandA(gm)#" (default formatting)
$sub = $d.Range($p, $p + 6)
$sub.Font.Color = 255  # w:color FF0000
$p = $p + 6
$p = $p + 27  # run 3: "opaq CandA(gm)#siq CandA(gm" (default formatting)
$sub = $d.Range($p, $p + 8)
$sub.Font.Color = 255  # w:color FF0000
$p = $p + 8
$p = $p + 4  # run 5: "a# |" (default formatting)
$pos = $paraTextEnd + 1
